# Auto commit at 2026-01-11  7:59:58.85
# Append two new rows (20 and 21) of day-data for 2026-01-10 (serial 46032)
# for station "四方坪站" and "高岭站" respectively, matching the style of
# the existing data rows, then move the active selection to J20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - 四方坪站
$ws.Cells.Item(20, 1).Value = 46032
$ws.Cells.Item(20, 2).Value = "四方坪站"
$ws.Cells.Item(20, 3).Value = 14688.8
$ws.Cells.Item(20, 4).Value = 10309.1
$ws.Cells.Item(20, 5).Value = 3244.76
$ws.Cells.Item(20, 6).Value = 648

# Row 21 - 高岭站
$ws.Cells.Item(21, 1).Value = 46032
$ws.Cells.Item(21, 2).Value = "高岭站"
$ws.Cells.Item(21, 3).Value = 4882.21
$ws.Cells.Item(21, 4).Value = 3838.73
$ws.Cells.Item(21, 5).Value = 1374.38
$ws.Cells.Item(21, 6).Value = 167

# Update the selected cell shown in the saved view, as in the source workbook
$ws.Range("J20").Select()
